$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns J (pixels) and K (rgb nums) ---
$ws.Range("J1").Value = "pixels"
$ws.Range("K1").Value = "rgb nums"

# --- Two new data rows (5 and 6) ---
$ws.Range("A5").Value = 25
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.1
$ws.Range("D5").Value = 500
$ws.Range("E5").Value = "siren"
$ws.Range("F5").Value = 771
$ws.Range("G5").Value = 0.0000027336
$ws.Range("G5").NumberFormat = "0.0000000000"
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 128

$ws.Range("A6").Value = 25
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.05
$ws.Range("D6").Value = 500
$ws.Range("E6").Value = "siren"
$ws.Range("F6").Value = 387
$ws.Range("G6").Value = 0.0000047325
$ws.Range("G6").NumberFormat = "0.0000000000"
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 64

# --- Formulas for pixels / rgb nums ---
# Row 2 gets its own standalone formula (entered first, independently).
$ws.Range("J2").Formula = "=216*318"
$ws.Range("K2").Formula = "=J2*3"

# Rows 3-6 are filled together as one block, producing a single shared
# formula group spanning J3:J6 / K3:K6.
$ws.Range("J3:J6").Formula = "=216*318"
$ws.Range("K3:K6").Formula = "=J3*3"

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$null = $ws.Range("J14").Select()
